$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1652

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H129").Value = 1111.2291
$ws.Range("I129").Value = 680.7
$ws.Range("J129").Value = 1224.5264
$ws.Range("K129").Value = 2042.1
$ws.Range("L129").Value = 3673.5792
$ws.Range("M129").Value = 2957.9
$ws.Range("N129").Value = -13673.5792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 426318.53
$ws.Range("I32").Value = 454586.16
$ws.Range("J32").Value = 120085.664
$ws.Range("K32").Value = 454586.16
$ws.Range("L32").Value = 120085.664
$ws.Range("M32").Value = -454299.16
$ws.Range("N32").Value = -120659.664

$ws.Range("H68").Value = 100000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 100000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101622

$ws.Range("H71").Value = 100000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 100000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -308112

$ws.Range("H122").Value = 1202
$ws.Range("I122").Value = 1033.4615
$ws.Range("J122").Value = 1749.75
$ws.Range("K122").Value = 3100.3845
$ws.Range("L122").Value = 5249.25
$ws.Range("M122").Value = -650.3844999999997
$ws.Range("N122").Value = -10149.25

$ws.Range("H123").Value = 28481.584
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 28481.584
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 28481.584
$ws.Range("N123").Value = -38281.584

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 9137.333000000001
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 9137.333000000001
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 9137.333000000001
$ws.Range("N6").Value = -9363.333000000001

$ws.Range("H13").Value = 67245
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 67245
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 67245
$ws.Range("N13").Value = -67581

$ws.Range("H50").Value = 67340
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 67340
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 67340
$ws.Range("N50").Value = -68488

$ws.Range("H97").Value = 1059.3334
$ws.Range("I97").Value = 1059.3334
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1059.3334
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -68.33339999999998
$ws.Range("N97").ClearContents()

$ws.Range("H105").Value = 3004.3635
$ws.Range("I105").Value = 2339.6
$ws.Range("J105").Value = 3558.3333
$ws.Range("K105").Value = 2339.6
$ws.Range("L105").Value = 3558.3333
$ws.Range("M105").Value = -592.5999999999999
$ws.Range("N105").Value = -7052.3333

$ws.Range("H106").Value = 64400
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 64400
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 64400
$ws.Range("N106").Value = -66924

$ws.Range("H115").Value = 79122.664
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 79122.664
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 79122.664
$ws.Range("N115").Value = -82256.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 373.33334
$ws.Range("I105").Value = 373.33334
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 373.33334
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1373.66666

$ws.Range("H119").Value = 100761
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 100761
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 100761
$ws.Range("N119").Value = -110437

$ws.Range("H141").Value = 193664.58
$ws.Range("I141").Value = 200000
$ws.Range("J141").Value = 192608.67
$ws.Range("K141").Value = 200000
$ws.Range("L141").Value = 192608.67
$ws.Range("M141").Value = -194820
$ws.Range("N141").Value = -202968.67

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1350.898
$ws.Range("I68").Value = 833.65216
$ws.Range("J68").Value = 1509.52
$ws.Range("K68").Value = 2500.95648
$ws.Range("L68").Value = 4528.559999999999
$ws.Range("M68").Value = -1689.95648
$ws.Range("N68").Value = -6150.559999999999

$ws.Range("H71").Value = 1350.898
$ws.Range("I71").Value = 833.65216
$ws.Range("J71").Value = 1509.52
$ws.Range("K71").Value = 7502.869439999999
$ws.Range("L71").Value = 13585.68
$ws.Range("M71").Value = -3446.869439999999
$ws.Range("N71").Value = -21697.68

$ws.Range("H76").Value = 2304.3333
$ws.Range("I76").Value = 913
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 2739
$ws.Range("L76").Value = 9000
$ws.Range("M76").Value = -2356
$ws.Range("N76").Value = -9766

$ws.Range("H79").Value = 2304.3333
$ws.Range("I79").Value = 913
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 2739
$ws.Range("L79").Value = 9000
$ws.Range("M79").Value = -1413
$ws.Range("N79").Value = -11652

$ws.Range("H113").Value = 853.2273
$ws.Range("I113").Value = 507.52
$ws.Range("J113").Value = 1308.1052
$ws.Range("K113").Value = 1522.56
$ws.Range("L113").Value = 3924.3156
$ws.Range("M113").Value = 647.4400000000001
$ws.Range("N113").Value = -8264.3156

$ws.Range("H132").Value = 3012.7163
$ws.Range("I132").Value = 2019.0857
$ws.Range("J132").Value = 3904.4358
$ws.Range("K132").Value = 18171.7713
$ws.Range("L132").Value = 35139.9222
$ws.Range("M132").Value = -15641.7713
$ws.Range("N132").Value = -40199.9222

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 23000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 23000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 23000
$ws.Range("N39").Value = -24064

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H70").Value = 5831.5312
$ws.Range("I70").Value = 5778.2607
$ws.Range("J70").Value = 5967.6665
$ws.Range("K70").Value = 5778.2607
$ws.Range("L70").Value = 5967.6665
$ws.Range("M70").Value = -5508.2607
$ws.Range("N70").Value = -6507.6665

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H73").Value = 5831.5312
$ws.Range("I73").Value = 5778.2607
$ws.Range("J73").Value = 5967.6665
$ws.Range("K73").Value = 5778.2607
$ws.Range("L73").Value = 5967.6665
$ws.Range("M73").Value = -4842.2607
$ws.Range("N73").Value = -7839.6665

$ws.Range("H86").Value = 67515.75
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 67515.75
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 67515.75
$ws.Range("N86").Value = -69887.75

$ws.Range("H89").Value = 67515.75
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 67515.75
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 202547.25
$ws.Range("N89").Value = -214403.25

$ws.Range("H99").Value = 11810.143
$ws.Range("I99").Value = 13167.75
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 13167.75
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = -10921.75
$ws.Range("N99").Value = -14492

$ws.Range("H126").Value = 2501
$ws.Range("I126").Value = 2501
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7503
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5033
$ws.Range("N126").ClearContents()

$ws.Range("H130").Value = 46593.332
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 46593.332
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 46593.332
$ws.Range("N130").Value = -56633.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 333336160
$ws.Range("I7").Value = 1000000000
$ws.Range("J7").Value = 4252.5
$ws.Range("K7").Value = 1000000000
$ws.Range("L7").Value = 4252.5
$ws.Range("M7").Value = -999999888
$ws.Range("N7").Value = -4476.5

$ws.Range("H40").Value = 142860860
$ws.Range("I40").Value = 333335170
$ws.Range("J40").Value = 5125
$ws.Range("K40").Value = 333335170
$ws.Range("L40").Value = 5125
$ws.Range("M40").Value = -333335034
$ws.Range("N40").Value = -5397

$ws.Range("H93").Value = 9000.643
$ws.Range("I93").Value = 21318.8
$ws.Range("J93").Value = 2157.2222
$ws.Range("K93").Value = 21318.8
$ws.Range("L93").Value = 2157.2222
$ws.Range("M93").Value = -20070.8
$ws.Range("N93").Value = -4653.2222

$ws.Range("H101").Value = 25472.4
$ws.Range("I101").Value = 5000
$ws.Range("J101").Value = 30590.5
$ws.Range("K101").Value = 5000
$ws.Range("L101").Value = 30590.5
$ws.Range("M101").Value = -1755
$ws.Range("N101").Value = -37080.5

$ws.Range("H122").Value = 6583.6665
$ws.Range("I122").Value = 5252
$ws.Range("J122").Value = 6850
$ws.Range("K122").Value = 15756
$ws.Range("L122").Value = 20550
$ws.Range("M122").Value = -13306
$ws.Range("N122").Value = -25450

$ws.Range("H126").Value = 333336160
$ws.Range("I126").Value = 1000000000
$ws.Range("J126").Value = 4252.5
$ws.Range("K126").Value = 3000000000
$ws.Range("L126").Value = 12757.5
$ws.Range("M126").Value = -2999997530
$ws.Range("N126").Value = -17697.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 10000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 10000
$ws.Range("N74").Value = -11872

$ws.Range("H77").Value = 10000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 30000
$ws.Range("N77").Value = -39360

$ws.Range("H123").Value = 24846.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 24846.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 24846.5
$ws.Range("N123").Value = -34646.5
